$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade rows to append (rows 7-9), matching columns:
# A=Principle, B=Start Principle, C=BuyPrice, D=SellPrice, E=IsShortSell,
# F=Price Change %, G=Date, H=Profitable

$data = @(
    @(9880.93, 9797.65, 18.84, 19,    $false, 0.85, 42613.766956018517, $true),
    @(9901.68, 9880.93, 18.93, 18.97, $false, 0.21, 42614.675405092596, $true),
    @(9991.7900000000009, 9901.68, 18.72, 18.89, $false, 0.91, 42615.752129629633, $true)
)

$rowIndex = 7
foreach ($rowData in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $rowData[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rowData[1]
    $ws.Cells.Item($rowIndex, 3).Value = $rowData[2]
    $ws.Cells.Item($rowIndex, 4).Value = $rowData[3]
    $ws.Cells.Item($rowIndex, 5).Value = $rowData[4]
    $ws.Cells.Item($rowIndex, 6).Value = $rowData[5]
    $ws.Cells.Item($rowIndex, 7).Value = $rowData[6]
    $ws.Cells.Item($rowIndex, 8).Value = $rowData[7]
    $rowIndex++
}

# Copy the date/time number format (style index referencing numFmtId 22)
# from the existing G6 cell onto the newly added G7:G9 cells, reusing the
# same style entry instead of creating a new custom number format.
$ws.Cells.Item(6, 7).Copy()
$ws.Range("G7:G9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

